$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '66.945.76'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '3.116.42'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '577.17'
$ws.Range('E5').Value = '  -0.69%  '
Set-TextValue 'D6' '171.53'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.113.59'
$ws.Range('E8').Value = '  +0.63%  '
Set-TextValue 'D9' '0.522'
$ws.Range('E9').Value = '  -0.80%  '
Set-TextValue 'D10' '6.50'
$ws.Range('E10').Value = '  -2.71%  '
Set-TextValue 'D11' '0.153'
$ws.Range('E11').Value = '  -1.39%  '
Set-TextValue 'D12' '0.483'
$ws.Range('E12').Value = '  +0.06%  '
Set-TextValue 'D13' '0.0000246'
$ws.Range('E13').Value = '  -1.59%  '
Set-TextValue 'D14' '37.19'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').Value = '3.628.95'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '66.903.91'
$ws.Range('E17').Value = '  -0.09%  '
Set-TextValue 'D18' '7.14'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').Value = '3.113.60'
$ws.Range('E19').Value = '  +0.58%  '
Set-TextValue 'D20' '16.37'
$ws.Range('E20').Value = '  -0.48%  '
Set-TextValue 'D21' '476.88'
$ws.Range('E21').Value = '  +1.44%  '
$ws.Range('E22').Value = '  -0.22%  '
Set-TextValue 'D23' '7.95'
$ws.Range('E23').Value = '  +5.32%  '
Set-TextValue 'D24' '13.48'
$ws.Range('E24').Value = '  +4.65%  '
Set-TextValue 'D25' '84.03'
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('E26').Value = '  -1.84%  '
Set-TextValue 'D27' '10.08'
$ws.Range('E27').Value = '  -0.61%  '
Set-TextValue 'D28' '1.00'
$ws.Range('E28').Value = '  +0.04%  '
Set-TextValue 'D29' '7.93'
$ws.Range('E29').Value = '  -1.34%  '
Set-TextValue 'D30' '2.40'
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').Value = '  +0.69%  '
Set-TextValue 'D33' '0.115'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '0.0₃0940'
$ws.Range('E34').Value = '  -7.74%  '
Set-TextValue 'D35' '0.999'
$ws.Range('E35').Value = '  -0.08%  '
Set-TextValue 'D36' '5.87'
$ws.Range('E36').Value = '  -0.97%  '
Set-TextValue 'D37' '0.977'
$ws.Range('E37').Value = '  -3.03%  '
Set-TextValue 'D38' '47.26'
$ws.Range('E38').Value = '  +1.13%  '
$ws.Range('E39').Value = '  -0.67%  '
Set-TextValue 'D40' '50.08'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('E41').Value = '  -2.04%  '
$ws.Range('E42').Value = '  -0.82%  '
Set-TextValue 'D43' '8.73'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '2.815.34'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0357'
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D46' '2.59'
$ws.Range('E46').Value = '  -10.26%  '
Set-TextValue 'D47' '381.73'
$ws.Range('E47').Value = '  -2.08%  '
Set-TextValue 'D48' '135.97'
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('E49').Value = '  +0.01%  '
Set-TextValue 'D50' '24.84'
$ws.Range('E50').Value = '  +0.19%  '
Set-TextValue 'D51' '2.20'
$ws.Range('E51').Value = '  -1.94%  '

Write-Host "Applied all changes"
